# This workbook was authored as a "Login" sheet (username/password demo data)
# with a Russian-localized Excel ("Лист1" sheet name, custom column width,
# custom row height, cell G9 selected). The commit re-saves the file with an
# older/plain Excel profile: the sheet is renamed to the default "Sheet1",
# the stray column-width/row-height customizations are cleared back to
# Excel's stock defaults, and the active selection becomes B2.
#
# Read all the existing cell values up front (handles can go stale once a
# new sheet is added), then build a brand-new, unformatted worksheet,
# transfer the values across, drop the old sheet, and rename the new one.

$wb = $excel.ActiveWorkbook
$old = $wb.Worksheets.Item(1)

$a1 = $old.Range("A1").Value2
$b1 = $old.Range("B1").Value2
$a2 = $old.Range("A2").Value2
$b2 = $old.Range("B2").Value2

$new = $wb.Worksheets.Add()

$new.Range("A1").Value = $a1
$new.Range("B1").Value = $b1
$new.Range("A2").Value = $a2
$new.Range("B2").Value = $b2

# Re-fetch the original sheet by its (now shifted) index and remove it so
# only the freshly-created, default-formatted sheet remains.
$oldRef = $wb.Worksheets.Item(2)
[void]$oldRef.Delete()

$new.Name = "Sheet1"
[void]$new.Range("B2").Select()
